$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 349, pushing existing rows 349:375 down to 350:376
$ws.Rows.Item(349).Insert()

# Populate the new row 349 with this week's record for
# Feria Lagunitas de Puerto Montt / Zapallo italiano.
# Columns A,B,C,E,F,G,H,I,R stay constant for this whole block, so copy
# them straight from the (now shifted) row 350, which still holds the
# previous top record. (Read via Value2 - .Value getter is unreliable
# in this runtime.)
$ws.Range("A349").Value = $ws.Range("A350").Value2
$ws.Range("B349").Value = $ws.Range("B350").Value2
$ws.Range("C349").Value = $ws.Range("C350").Value2
$ws.Range("D349").Value = 45013
$ws.Range("E349").Value = $ws.Range("E350").Value2
$ws.Range("F349").Value = $ws.Range("F350").Value2
$ws.Range("G349").Value = $ws.Range("G350").Value2
$ws.Range("H349").Value = $ws.Range("H350").Value2
$ws.Range("I349").Value = $ws.Range("I350").Value2
$ws.Range("J349").Value = 250
$ws.Range("K349").Value = 11000
$ws.Range("L349").Value = 11000
$ws.Range("M349").Value = 11000
$ws.Range("N349").Value = "$/caja 50 unidades"
$ws.Range("O349").Value = "Región Metropolitana"
$ws.Range("P349").Value = 220
$ws.Range("Q349").Value = 50
$ws.Range("R349").Value = $ws.Range("R350").Value2
